$d = $word.ActiveDocument

# --- Step 1: rewrite paragraphs 2-6 (Figure 1 .. Figure 5 legends) via InsertXML ---
# (done BEFORE deleting the trailing scratch paragraphs, so paragraph 6 is never the
#  last paragraph in the body while being rewritten -- InsertXML on the very last
#  paragraph of a body leaves a stray empty paragraph behind.)
$xml2 = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Figure 1</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">RTL in relation to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">(A)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">age and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">(B)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>Δ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">Age in Seychelles warblers. Lines and shaded areas are fitted values and 95% confidence limits from linear regressions.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(2).Range.InsertXML($xml2)

$xml3 = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Figure 2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Cohort effects on early-life RTL in Seychelles warblers.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">A</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">temporal variation RTL in major and minor cohorts,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">temporal variation in adult population size, and</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">C</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">early-life RTL in relation to cohort-level variation in adult population size. Lines and shaded areas in</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">C</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">are fitted values and 95% confidence limits from general linear models.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(3).Range.InsertXML($xml3)

$xml4 = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Figure 3</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Factors affecting early-life RTL in Seychelles warblers.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">A</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Model averaged estimates and 95% confidence intervals for explanatory terms used in a linear mixed model with RTL as the response variable (see main text for details; TQ = territory quality).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B,C</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Early-life RTL and age in relation to tarsus length (</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B</w:t></w:r><w:r><w:t xml:space="preserve">) and sex (</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">C</w:t></w:r><w:r><w:t xml:space="preserve">). Lines and shaded areas in</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">are fitted values and 95% confidence limits from a general linear model.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(4).Range.InsertXML($xml4)

$xml5 = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Figure</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Longitudinal telomere dynamics in juvenile Seychelles warblers.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">A</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Juvenile telomere length in relation to adult telomere length measured in the same individual.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>Δ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">RTL in relation to the time between sampling events. Lines and shaded areas represent fitted values and 95% confidence limits from general linear models. Raw data have been ommited from</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">for the sake of clarity, but are plotted in Figure S1.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(5).Range.InsertXML($xml5)

$xml6 = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Figure 5</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Factors affecting</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>Δ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">RTL in juvenile Seychelles warblers.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">A</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Model averaged estimates and 95% confidence intervals for explanatory terms used in a linear mixed model with</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>Δ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">RTL as the response variable (see main text for details; TQ = territory quality).</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">B</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:rPr><m:sty m:val="p"/></m:rPr><m:t>Δ</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">RTL in relation to tarsus length and sex.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(6).Range.InsertXML($xml6)

# --- Step 2: remove the trailing scratch/code paragraphs (old paragraphs 7-14) ---
$delStart = $d.Paragraphs.Item(7).Range.Start
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$delEnd = $lastPara.Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : $($p.Range.Text)"
}
